$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de report sheets, as produced by a fresh
# "Generate Report for Handback" run.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 17:17:36"
$wsZhCn.Range("H2").Value = "2016-03-18 17:18:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 17:17:44"
$wsDeDe.Range("H2").Value = "2016-03-18 17:18:32"
